$wb = $excel.ActiveWorkbook

# --- "built to spill" sheet: duplicate columns A:B into C:D, then correct the filenames in C ---
$ws7 = $wb.Worksheets.Item("built to spill")
$ws7.Range("A1:B58").Copy()
$ws7.Range("C1").PasteSpecial()
$ws7.Range("A1").Select()

# Retype the filenames that were cleaned up / zero-padded (creates the needed new entries
# in the shared string table in the same order the source workbook has them)
$ws7.Range("C31").Value = '08 bts-center-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C32").Value = '08 bts-fly-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C33").Value = '08 bts-goingagainst-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C34").Value = '08 bts-inthemorning-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C35").Value = '09 bts-madeupdreams-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C36").Value = '09 bts-sidewalk-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C37").Value = '09 bts-soundcheck-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C38").Value = '09 bts-theplan-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C39").Value = '09 bts-untitled-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C40").Value = '09 bts-youwereright-BuiltToSpill-Geogetown-20130713.mp3'
$ws7.Range("C42").Value = '02 In The Morning-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C43").Value = '03 Center of the Universe-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C44").Value = '04 The Plan-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C41").Value = '01 Going Against Your Mind-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C45").Value = '05 Planting Seeds-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C46").Value = '06 Kicked It In the Sun-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C47").Value = '07 Get a Life-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C48").Value = '08 Sidewalk-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C49").Value = '09 Here - Pavement-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C50").Value = '10 You were right-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C51").Value = '11 Conventional Wisdom-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C52").Value = '12 Heart - Doug Martsch-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C53").Value = '13 Carry the Zero-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C54").Value = '14 Instrumental - dont know-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C55").Value = '15 I would hurt a fly-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C56").Value = '16 Age of Consent - New Order-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C57").Value = '17 How Soon is Now - The Smiths-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C58").Value = '18 Car-BuiltToSpill-LiveatNeumos20131227.mp3'
$ws7.Range("C1").Value = '01%20Hazy2013-02-03-15-02-09-0679196-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C2").Value = '02%20The%20Source2013-02-03-15-02-10-4095712-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C3").Value = '03%20Reasons2013-02-03-15-02-11-6108174-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C4").Value = '05%20Stab2013-02-03-15-02-13-9977092-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C5").Value = '06%20Strange2013-02-03-15-02-15-2301566-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C6").Value = '07%20Get%20a%20Life2013-02-03-15-02-15-5421686-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C7").Value = '08%20Made%20Up%20Dreams2013-02-03-15-02-16-8214178-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C8").Value = '09%20Dont%20Fear%20the%20Reaper2013-02-03-15-17-23-8247012-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C9").Value = '10%20Goin%20Against%20Your%20Mind2013-02-03-15-17-25-1663528-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C10").Value = '11%20Kicked%20It%20In%20The%20Sun2013-02-03-15-17-25-9151816-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C11").Value = '12%20Center%20Of%20The%20Universe2013-02-03-15-17-27-2100314-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C12").Value = '13%20Heart2013-02-03-15-17-27-3192356-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C13").Value = '14%20Pat2013-02-03-15-17-28-4736800-BuiltToSpill-Bellingham-20130202.mp3'
$ws7.Range("C14").Value = '15%20How%20Soon%20Is%20Now%20-%20Unknown2013-02-03-15-17-29-5657220-BuiltToSpill-Bellingham-20130202.mp3'

$ws7.Columns.Item(3).ColumnWidth = 77.15
$ws7.Range("C1:C14").Select()

# --- "OPTIONaudio" sheet: point the Bellingham rows (221-234) at the corrected filenames ---
$ws1 = $wb.Worksheets.Item("OPTIONaudio")
$ws1.Range("B221").Value = '01%20Hazy2013-02-03-15-02-09-0679196-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B222").Value = '02%20The%20Source2013-02-03-15-02-10-4095712-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B223").Value = '03%20Reasons2013-02-03-15-02-11-6108174-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B224").Value = '05%20Stab2013-02-03-15-02-13-9977092-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B225").Value = '06%20Strange2013-02-03-15-02-15-2301566-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B226").Value = '07%20Get%20a%20Life2013-02-03-15-02-15-5421686-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B227").Value = '08%20Made%20Up%20Dreams2013-02-03-15-02-16-8214178-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B228").Value = '09%20Dont%20Fear%20the%20Reaper2013-02-03-15-17-23-8247012-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B229").Value = '10%20Goin%20Against%20Your%20Mind2013-02-03-15-17-25-1663528-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B230").Value = '11%20Kicked%20It%20In%20The%20Sun2013-02-03-15-17-25-9151816-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B231").Value = '12%20Center%20Of%20The%20Universe2013-02-03-15-17-27-2100314-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B232").Value = '13%20Heart2013-02-03-15-17-27-3192356-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B233").Value = '14%20Pat2013-02-03-15-17-28-4736800-BuiltToSpill-Bellingham-20130202.mp3'
$ws1.Range("B234").Value = '15%20How%20Soon%20Is%20Now%20-%20Unknown2013-02-03-15-17-29-5657220-BuiltToSpill-Bellingham-20130202.mp3'

$ws1.Range("D234").Select()
